# Updated cryptos list on Sat Jul 20 07:27:55 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced to Text
# (via Text number format) so they stay literal strings like "1.00" / "0.0313"
# instead of being auto-coerced to numeric cells, matching the source data which
# stores every Price/Volume column as inline text.
$textCells = @("D5", "D6", "D12", "D15", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D33", "D36", "D37", "D39", "D40", "D41", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.725.53"
$ws.Range("E2").Value = "  +4.27%  "
$ws.Range("D3").Value = "3.500.94"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "592.45"
$ws.Range("E5").Value = "  +3.47%  "
$ws.Range("D6").Value = "169.31"
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +8.41%  "
$ws.Range("D9").Value = "3.499.91"
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("E10").Value = "  +7.11%  "
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("D13").Value = "4.105.61"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "28.31"
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").Value = "66.707.67"
$ws.Range("E17").Value = "  +4.15%  "
$ws.Range("D18").Value = "3.488.39"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").Value = "14.21"
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("D21").Value = "393.26"
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("D22").Value = "7.99"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").Value = "73.45"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "0.538"
$ws.Range("E25").Value = "  +4.17%  "
$ws.Range("E26").Value = "  +5.70%  "
$ws.Range("D27").Value = "10.26"
$ws.Range("E27").Value = "  +7.21%  "
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "6.41"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("D33").Value = "23.66"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("E34").Value = "  +4.50%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "1.64"
$ws.Range("E36").Value = "  +8.97%  "
$ws.Range("D37").Value = "162.74"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").Value = "  +5.89%  "
$ws.Range("D40").Value = "27.71"
$ws.Range("E40").Value = "  +5.35%  "
$ws.Range("D41").Value = "6.79"
$ws.Range("E41").Value = "  +4.83%  "
$ws.Range("E42").Value = "  +5.96%  "
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").Value = "26.60"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("D45").Value = "2.795.53"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "43.23"
$ws.Range("D47").Value = "0.0313"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").Value = "353.08"
$ws.Range("E49").Value = "  +7.29%  "
$ws.Range("E50").Value = "  +5.49%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "0.885"
$ws.Range("E51").Value = "  +8.47%  "
